$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.576.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "'2.314.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'318.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.28%  "
$ws.Range("D6").Value = "'104.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("D7").Value = "'0.635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.613"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("D10").Value = "'39.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").Value = "'0.0912"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "'8.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "'0.973"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "'15.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").Value = "'2.663.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "'2.311.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").Value = "'42.507.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").Value = "'7.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").Value = "'3.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.24%  "
$ws.Range("D22").Value = "'73.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("D23").Value = "'282.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.33%  "
$ws.Range("D24").Value = "'11.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +20.73%  "
$ws.Range("D25").Value = "'2.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "'10.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("E28").Value = "  +6.26%  "
$ws.Range("D29").Value = "'23.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "'36.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.24%  "
$ws.Range("D31").Value = "'165.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "'0.0879"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("D33").Value = "'5.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("E34").Value = "  +5.19%  "
$ws.Range("E35").Value = "  -7.61%  "
$ws.Range("D36").Value = "'0.116"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("D37").Value = "'4.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.02%  "
$ws.Range("D38").Value = "'0.0352"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").Value = "'3.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("D40").Value = "'2.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.42%  "
$ws.Range("D41").Value = "'101.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.23%  "
$ws.Range("D42").Value = "'1.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("D43").Value = "'69.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("D44").Value = "'0.228"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.50%  "
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "'12.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'113.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("D48").Value = "'78.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.01%  "
$ws.Range("D49").Value = "'8.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").Value = "'5.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("D51").Value = "'1.612.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.13%  "
